# calorimetry : input and output consistency : done
#
# The "enthalpies_calculated" sheet and the "input_enthalpies" sheet both
# used to carry rows for every species (H, L, M, HL, H2L, ML, HML) even
# though only the complexes (HL, H2L, ML, HML) actually have enthalpy
# data. This drops the H/L/M placeholder rows (which were all zero) and
# fixes the sign on the ML / HML calculated enthalpy values so the
# calculated-output sheet is consistent with the input sheet.

$wb = $excel.ActiveWorkbook

# --- enthalpies_calculated --------------------------------------------
$ws = $wb.Worksheets.Item("enthalpies_calculated")

# Remove the H, L, M rows (rows 2-4); this shifts HL/H2L/ML/HML up so
# they become rows 2-5.
$ws.Rows("2:4").Delete()

# The ML and HML calculated values were stored with the wrong sign -
# correct them.
$ws.Range("B4").Value = -9.47878966930947
$ws.Range("B5").Value = -11.2246759055297

# --- input_enthalpies ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("input_enthalpies")

# Same cleanup: drop the H, L, M rows so only HL and H2L remain.
$ws2.Rows("2:4").Delete()
